$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert row 8 (inherits formatting from row 7 above, matching Excel's
# native "insert row" behaviour) and populate the new EV/FR00 entry.
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "FR00"
$ws.Range("B8").Value = "Distributed Energy"
$ws.Range("C8").Value = 2030
$ws.Range("D8").Value = 150000

# Insert row 9 (inherits formatting from row 8 above) and populate the
# second new EV/FR00 entry.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "FR00"
$ws.Range("B9").Value = "Distributed Energy"
$ws.Range("C9").Value = 2040
$ws.Range("D9").Value = 300000

# Match the author's final selection/active cell left in the file.
$ws.Range("L5").Select() | Out-Null
